$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.880.16"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "3.324.64"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'584.24"
$ws.Range("E5").Value = "  +3.93%  "
$ws.Range("D6").Value = "'182.80"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").Value = "3.321.82"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "'0.179"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'46.40"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").Value = "'635.38"
$ws.Range("E14").Value = "  +8.10%  "
$ws.Range("D15").Value = "3.858.05"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "67.981.15"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "3.324.60"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'17.70"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'10.93"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'0.903"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'17.68"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'97.05"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "'4.00"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'2.78"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").Value = "'32.57"
$ws.Range("E29").Value = "  +6.64%  "
$ws.Range("D30").Value = "'8.63"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("D31").Value = "'6.73"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("D32").Value = "'594.47"
$ws.Range("E32").Value = "  +6.00%  "
$ws.Range("D33").Value = "3.954.80"
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("D34").Value = "'10.95"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "'0.105"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").Value = "'3.53"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'55.67"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("D41").Value = "'2.69"
$ws.Range("E41").Value = "  +4.85%  "
$ws.Range("D42").Value = "'32.63"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "0.0₃0686"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "'0.0414"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  +12.80%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'130.46"
$ws.Range("E51").Value = "  +2.06%  "
